$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("optimization_parameters")
$ws.Rows.Item(16).Delete()
Write-Output $ws.Range("A16").Value()
Write-Output $ws.Range("B16").Value()
Write-Output $ws.Range("C16").Value()
Write-Output $ws.Range("A17").Value()
